$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 4).Value = 44315

# Row 3
$ws.Cells.Item(3, 4).Value = 44309
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 45

# Row 4
$ws.Cells.Item(4, 4).Value = 44302
$ws.Cells.Item(4, 13).Value = 45
$ws.Cells.Item(4, 14).Value = 10000
$ws.Cells.Item(4, 15).Value = 10000
$ws.Cells.Item(4, 16).Value = 10000
$ws.Cells.Item(4, 19).Value = 1000

# Row 5
$ws.Cells.Item(5, 4).Value = 44307
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 40
$ws.Cells.Item(5, 14).Value = 10000
$ws.Cells.Item(5, 15).Value = 10000
$ws.Cells.Item(5, 16).Value = 10000
$ws.Cells.Item(5, 19).Value = 1000

# Row 6
$ws.Cells.Item(6, 4).Value = 44301
$ws.Cells.Item(6, 13).Value = 45
$ws.Cells.Item(6, 14).Value = 10000
$ws.Cells.Item(6, 15).Value = 10000
$ws.Cells.Item(6, 16).Value = 10000
$ws.Cells.Item(6, 19).Value = 1000

# Row 7
$ws.Cells.Item(7, 4).Value = 44323
$ws.Cells.Item(7, 13).Value = 60

# Row 8
$ws.Cells.Item(8, 4).Value = 44323
$ws.Cells.Item(8, 12).Value = 'Segunda'
$ws.Cells.Item(8, 13).Value = 50
$ws.Cells.Item(8, 14).Value = 9000
$ws.Cells.Item(8, 15).Value = 9000
$ws.Cells.Item(8, 16).Value = 9000
$ws.Cells.Item(8, 19).Value = 900

# Row 9
$ws.Cells.Item(9, 4).Value = 44326
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 65
$ws.Cells.Item(9, 14).Value = 10000
$ws.Cells.Item(9, 15).Value = 10000
$ws.Cells.Item(9, 16).Value = 10000
$ws.Cells.Item(9, 19).Value = 1000

# Row 10
$ws.Cells.Item(10, 4).Value = 44326
$ws.Cells.Item(10, 12).Value = 'Segunda'
$ws.Cells.Item(10, 13).Value = 67
$ws.Cells.Item(10, 14).Value = 8000
$ws.Cells.Item(10, 15).Value = 8000
$ws.Cells.Item(10, 16).Value = 8000
$ws.Cells.Item(10, 19).Value = 800

# Row 11
$ws.Cells.Item(11, 4).Value = 44322
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 56
$ws.Cells.Item(11, 14).Value = 10000
$ws.Cells.Item(11, 15).Value = 10000
$ws.Cells.Item(11, 16).Value = 10000
$ws.Cells.Item(11, 19).Value = 1000

# Row 12
$ws.Cells.Item(12, 4).Value = 44322
$ws.Cells.Item(12, 12).Value = 'Segunda'
$ws.Cells.Item(12, 13).Value = 40
$ws.Cells.Item(12, 14).Value = 8000
$ws.Cells.Item(12, 15).Value = 8000
$ws.Cells.Item(12, 16).Value = 8000
$ws.Cells.Item(12, 19).Value = 800

# Row 15
$ws.Cells.Item(15, 4).Value = 44308

# Row 16
$ws.Cells.Item(16, 4).Value = 44308
$ws.Cells.Item(16, 12).Value = 'Segunda'
$ws.Cells.Item(16, 13).Value = 48
$ws.Cells.Item(16, 14).Value = 8000
$ws.Cells.Item(16, 15).Value = 8000
$ws.Cells.Item(16, 16).Value = 8000
$ws.Cells.Item(16, 19).Value = 800

# Row 17
$ws.Cells.Item(17, 4).Value = 44312
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 48
$ws.Cells.Item(17, 18).Value = 'Provincia de Quillota'

# Row 18
$ws.Cells.Item(18, 4).Value = 44333
$ws.Cells.Item(18, 12).Value = 'Especial'
$ws.Cells.Item(18, 13).Value = 58
$ws.Cells.Item(18, 14).Value = 10000
$ws.Cells.Item(18, 15).Value = 10000
$ws.Cells.Item(18, 16).Value = 10000
$ws.Cells.Item(18, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(18, 19).Value = 1000

# Row 19
$ws.Cells.Item(19, 4).Value = 44333
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 65
$ws.Cells.Item(19, 14).Value = 9000
$ws.Cells.Item(19, 15).Value = 9000
$ws.Cells.Item(19, 16).Value = 9000
$ws.Cells.Item(19, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(19, 19).Value = 900

# Row 20
$ws.Cells.Item(20, 4).Value = 44333
$ws.Cells.Item(20, 12).Value = 'Segunda'
$ws.Cells.Item(20, 13).Value = 60
$ws.Cells.Item(20, 14).Value = 8000
$ws.Cells.Item(20, 15).Value = 8000
$ws.Cells.Item(20, 16).Value = 8000
$ws.Cells.Item(20, 19).Value = 800

# Row 21
$ws.Cells.Item(21, 4).Value = 44306
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 45
$ws.Cells.Item(21, 14).Value = 10000
$ws.Cells.Item(21, 15).Value = 10000
$ws.Cells.Item(21, 16).Value = 10000
$ws.Cells.Item(21, 19).Value = 1000

# Row 22
$ws.Cells.Item(22, 4).Value = 44328
$ws.Cells.Item(22, 14).Value = 8000
$ws.Cells.Item(22, 15).Value = 8000
$ws.Cells.Item(22, 16).Value = 8000
$ws.Cells.Item(22, 19).Value = 800

# Row 23
$ws.Cells.Item(23, 4).Value = 44328
$ws.Cells.Item(23, 14).Value = 7000
$ws.Cells.Item(23, 15).Value = 7000
$ws.Cells.Item(23, 16).Value = 7000
$ws.Cells.Item(23, 19).Value = 700

# Row 24
$ws.Cells.Item(24, 4).Value = 44321
$ws.Cells.Item(24, 13).Value = 58
$ws.Cells.Item(24, 14).Value = 9000
$ws.Cells.Item(24, 15).Value = 9000
$ws.Cells.Item(24, 16).Value = 9000
$ws.Cells.Item(24, 19).Value = 900

# Row 25
$ws.Cells.Item(25, 4).Value = 44329
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 56
$ws.Cells.Item(25, 18).Value = 'Región Metropolitana'

# Row 26
$ws.Cells.Item(26, 4).Value = 44329
$ws.Cells.Item(26, 12).Value = 'Segunda'
$ws.Cells.Item(26, 13).Value = 50
$ws.Cells.Item(26, 14).Value = 8000
$ws.Cells.Item(26, 15).Value = 8000
$ws.Cells.Item(26, 16).Value = 8000
$ws.Cells.Item(26, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(26, 19).Value = 800

# Row 27
$ws.Cells.Item(27, 4).Value = 44319
$ws.Cells.Item(27, 13).Value = 68
$ws.Cells.Item(27, 14).Value = 10000
$ws.Cells.Item(27, 15).Value = 10000
$ws.Cells.Item(27, 16).Value = 10000
$ws.Cells.Item(27, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(27, 19).Value = 1000

# Row 28
$ws.Cells.Item(28, 4).Value = 44319
$ws.Cells.Item(28, 13).Value = 57
$ws.Cells.Item(28, 18).Value = 'Provincia de Quillota'

# Row 29
$ws.Cells.Item(29, 4).Value = 44343
$ws.Cells.Item(29, 12).Value = 'Especial'
$ws.Cells.Item(29, 13).Value = 47
$ws.Cells.Item(29, 18).Value = 'Región Metropolitana'

# Row 30
$ws.Cells.Item(30, 4).Value = 44343
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 50
$ws.Cells.Item(30, 14).Value = 9000
$ws.Cells.Item(30, 15).Value = 9000
$ws.Cells.Item(30, 16).Value = 9000
$ws.Cells.Item(30, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(30, 19).Value = 900

# Row 31
$ws.Cells.Item(31, 4).Value = 44343
$ws.Cells.Item(31, 12).Value = 'Segunda'
$ws.Cells.Item(31, 13).Value = 58
$ws.Cells.Item(31, 14).Value = 8000
$ws.Cells.Item(31, 15).Value = 8000
$ws.Cells.Item(31, 16).Value = 8000
$ws.Cells.Item(31, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(31, 19).Value = 800

# Row 32
$ws.Cells.Item(32, 4).Value = 44314
$ws.Cells.Item(32, 13).Value = 47
$ws.Cells.Item(32, 14).Value = 9000
$ws.Cells.Item(32, 15).Value = 9000
$ws.Cells.Item(32, 16).Value = 9000
